$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting of the two preceding rows (65 -> 67, 66 -> 68) so the
# new rows inherit the same alternating style used throughout the table,
# then fill in the new data values.
$ws.Range("A65:J65").Copy()
$ws.Range("A67:J67").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A66:J66").Copy()
$ws.Range("A68:J68").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$row67 = @(43967, 69363, 511, 1466, 1, 26, 5, 1, 104, 1)
$row68 = @(43968, 69842, 479, 1466, 0, 25, 5, 1, 104, 0)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(67, $i + 1).Value = $row67[$i]
    $ws.Cells.Item(68, $i + 1).Value = $row68[$i]
}

# Extend the table to include the two new rows.
$tbl = $ws.ListObjects.Item("Tabela1")
$tbl.Resize($ws.Range("A1:J68"))

# Update dimension-driving selection/view state to match the new extent.
$ws.Range("A68:J68").Select()
